$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptos list on Thu May 11 06:50:20 UTC 2023 with GitHub Actions
# Refresh Price (D) and Volume(1h) (E) columns for each coin row.
# Price cells are forced to Text format before the write so values like
# "312.30" / "1.003" are stored verbatim (not coerced to numbers), then
# the style is reset back to Normal so no stray formatting is introduced.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.527.43"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.832.34"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.68%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "312.30"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.18%  "
$ws.Range("E6").Value = "  -0.06%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4254"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.75%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3648"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.33%  "
$ws.Range("E9").Value = "  -0.83%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8623"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.74%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.61"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.35%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.806.21"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.02%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.398"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.90%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.505"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.19%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.06936"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.27%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.003"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.10%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "80.46"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008860"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.25%  "
$ws.Range("E19").Value = "  -0.06%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.38"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.03%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "27.516.08"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.50%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.147"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.28%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.82"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +4.94%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.055.26"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.01%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.989"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.10%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "154.47"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.84%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.79"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.35%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.106"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.09%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "114.33"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.41%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.815"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.00%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08836"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.54%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.983"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.61%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7447"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.09%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.526"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.131"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.73%  "
$ws.Range("E36").Value = "  -0.01%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.089"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.48%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05300"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.61%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.797"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.26%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5068"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.06%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1654"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.62%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.480"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.62%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.308"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.67%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.42"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.63%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.06471"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.02%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "105.22"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.61%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.4672"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.45%  "
$ws.Range("E49").Value = "  -0.09%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.612"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.34%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "63.41"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.57%  "
